$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching the formatting of the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# New data cells for the "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
